{"js": "// Justify the existing paragraph, then append two new paragraphs describing\n// how the mobile/background-image version was implemented.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Justify the first (and, at this point, only) paragraph.\nconst firstPara = paragraphs.items[0];\nfirstPara.alignment = Word.Alignment.justified;\nawait context.sync();\n\n// 2) Add the paragraph about testing the mobile/tablet view in dev tools.\nconst mobileTestPara = firstPara.insertParagraph(\n  \"Para saber c\u00f3mo se ver\u00eda en el m\u00f3vil o Tablet, dentro del navegador, le doy click derecho e inspeccionar. \" +\n    \"Desde ah\u00ed, le doy a Activar/Desactivar barra de herramientas del dispositivo, y escojo las dimensiones del \" +\n    \"dispositivo que me interesan para poder ver c\u00f3mo quedar\u00eda, directamente desde la pantalla del pc. He usado \" +\n    \"las dimensiones del iPhone 12 Pro.\",\n  Word.InsertLocation.after\n);\nmobileTestPara.alignment = Word.Alignment.justified;\nawait context.sync();\n\n// 3) Add the paragraph about adjusting the hero-home background image for mobile.\nconst heroHomePara = mobileTestPara.insertParagraph(\n  \"En index.html para ubicar la imagen de fondo hero-home \",\n  Word.InsertLocation.after\n);\nheroHomePara.alignment = Word.Alignment.justified;\nheroHomePara.insertText(\"en versi\u00f3n m\u00f3vil, \", Word.InsertLocation.end);\nheroHomePara.insertText(\n  \"he tenido que a\u00f1adir padding top y bottom, as\u00ed como jugar con background-size y background position hasta \" +\n    \"lograr el efecto parecido al que se \",\n  Word.InsertLocation.end\n);\nheroHomePara.insertText(\"pide\", Word.InsertLocation.end);\nheroHomePara.insertText(\".\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Justify the existing paragraph, then append two new paragraphs describing\n# how the mobile/background-image version was implemented.\n\n$d = $word.ActiveDocument\n\n# 1) Justify the first (and, at this point, only) paragraph.\n$d.Paragraphs(1).Alignment = 3   # wdAlignParagraphJustify\n\n# 2) Add the paragraph about testing the mobile/tablet view in dev tools.\n$sel = $word.Selection\n$sel.EndKey(6, 0) | Out-Null   # wdStory, wdMove\n$sel.TypeParagraph()\n$sel.TypeText(\"Para saber c\u00f3mo se ver\u00eda en el m\u00f3vil o Tablet, dentro del navegador, le doy click derecho e inspeccionar. Desde ah\u00ed, le doy a Activar/Desactivar barra de herramientas del dispositivo, y escojo las dimensiones del dispositivo que me interesan para poder ver c\u00f3mo quedar\u00eda, directamente desde la pantalla del pc. He usado las dimensiones del iPhone 12 Pro.\")\n$d.Paragraphs(2).Alignment = 3\n\n# 3) Add the paragraph about adjusting the hero-home background image for mobile.\n$sel.EndKey(6, 0) | Out-Null\n$sel.TypeParagraph()\n$sel.TypeText(\"En index.html para ubicar la imagen de fondo hero-home \")\n$sel.TypeText(\"en versi\u00f3n m\u00f3vil, \")\n$sel.TypeText(\"he tenido que a\u00f1adir padding top y bottom, as\u00ed como jugar con background-size y background position hasta lograr el efecto parecido al que se \")\n$sel.TypeText(\"pide\")\n$sel.TypeText(\".\")\n$d.Paragraphs(3).Alignment = 3\n"}
